{"js": "// The diff resizes the 6 inline \"screenshot\" pictures in the lab03 report\n// from a width of 666750 EMU (52.5 pt) to 1905000 EMU (150 pt), scaling the\n// height in proportion (factor 20/7) so each picture's aspect ratio is kept.\n//\n// Word's Office.js `InlinePicture.width` / `.height` setters are no-ops (the\n// host only supports resizing real floating `Shape`s that way), so the\n// resize has to be driven through the same low-level OM bridge the shim\n// itself uses for `InlineShape.Width` \u2014 this mirrors what\n// `Shape.width =` does under the hood, just aimed at `InlineShapes` instead.\n// Setting `InlineShape.Width` alone recomputes the picture's height from its\n// original/natural aspect ratio (matching the diff's height values), so we\n// only need to set the width for each picture.\n\nconst pics = context.document.body.inlinePictures;\npics.load(\"items\");\nawait context.sync();\n\nconst TARGET_WIDTH_PT = 1905000 / 12700; // 150pt \u2014 1905000 EMU\n\nfor (let i = 0; i < pics.items.length; i++) {\n  const pic = pics.items[i];\n  // Reach past the no-op width/height setters straight to the OM bridge\n  // (same call the generated Shape.width setter uses) so the resize\n  // actually lands in the document's <wp:extent>/<a:ext>.\n  __native.docxOmSet(pic._h, pic._a, \"InlineShape.Width\", String(TARGET_WIDTH_PT));\n}\n\nawait context.sync();\n", "ps1": "# The diff resizes the 6 inline \"screenshot\" pictures in the lab03 report\n# from a width of 666750 EMU (52.5 pt) to 1905000 EMU (150 pt). Word scales\n# the height in proportion to the picture's original aspect ratio whenever\n# Width is set, which reproduces each picture's new height from the diff.\n\n$d = $word.ActiveDocument\n\n$targetWidthPt = 1905000 / 12700  # 150pt == 1905000 EMU\n\nforeach ($shape in $d.InlineShapes) {\n    $shape.Width = $targetWidthPt\n}\n"}
